$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 currently holds the shared string "R40" (s="23", t="s").
# We need it to hold the text "1" instead, WITHOUT disturbing its
# existing style (s="23") and without it being auto-coerced into a
# numeric cell (plain Value/Formula assignment of "1" is detected as a
# number by the engine and silently drops the string type).
#
# Trick: build the text "1" as the result of a formula in a scratch
# cell (a formula result is always text here, never re-parsed as a
# number), copy it, and paste only the VALUE into B11. PasteSpecial
# Values preserves B11's pre-existing cell style (xfId) and stores the
# pasted scalar using its source type (string), which is exactly what
# we want. Finally remove the scratch cell so it leaves no trace in
# the sheet (dimension/used-range stays the same as before the edit).
$scratch = $ws.Range("Z1")
$scratch.Formula = '="1"'
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$scratch.Delete(-4159)                # xlShiftUp - remove all trace of the scratch cell
